$wb = $excel.ActiveWorkbook

# --- Sheet "DPLKKPS144-001" (row 2 = "Kembalikan ke Register" scenario) ---
# Only the No. Register value changes (...028 -> ...027); the Saldo Nominal
# Final balances (O2/P2) stay the same.
$ws1 = $wb.Worksheets.Item("DPLKKPS144-001")

$ws1.Range("N2").Value = "M11220800000027"

$prep1 = @"
Username : 30711;
Password : bni1234;
Role : 09 - Penyelia Settlement;
No. Register : M11220800000027;
Saldo Nominal Final - Saldo Awal Iuran Pribadi : 308.000,00;
Saldo Nominal Final - Saldo Awal Iuran Perusahaan : 1.292.000,00;
Saldo Nominal Final - Saldo Awal Iuran Sukarela : 0,00;
Saldo Nominal Final - Saldo Awal Pengalihan Iuran Karyawan : 0,00;
Saldo Nominal Final - Saldo Awal Pengalihan Iuran Perusahaan : 0,00;
Status Verifikasi : 0 : Kembalikan ke Register;
Keterangan Verifikasi : KEP.TRX.445 Data dikembalikan
"@
$ws1.Range("F2").Value = $prep1

# --- Sheet "DPLKKPS144-002" (row 2 = "Setuju" scenario) ---
# No. Register changes (...028 -> ...011) AND the Saldo Nominal Final
# balances change (Iuran Pribadi / Iuran Perusahaan).
$ws2 = $wb.Worksheets.Item("DPLKKPS144-002")

$ws2.Range("N2").Value = "M11220800000011"
$ws2.Range("O2").Value = 2000000
$ws2.Range("P2").Value = 1600000

$prep2 = @"
Username : 30711;
Password : bni1234;
Role : 09 - Penyelia Settlement;
No. Register : M11220800000011;
Saldo Nominal Final - Saldo Awal Iuran Pribadi : 2.000.000,00;
Saldo Nominal Final - Saldo Awal Iuran Perusahaan : 1.600.000,00;
Saldo Nominal Final - Saldo Awal Iuran Sukarela : 0,00;
Saldo Nominal Final - Saldo Awal Pengalihan Iuran Karyawan : 0,00;
Saldo Nominal Final - Saldo Awal Pengalihan Iuran Perusahaan : 0,00;
Status Verifikasi : 1 : Setuju;
Keterangan Verifikasi : KEP.TRX.445 Disetujui
"@
$ws2.Range("F2").Value = $prep2

# --- View state: sheet1 loses focus/selection, sheet2 becomes the active tab ---
$ws1.Activate() | Out-Null
$excel.ActiveWindow.ScrollColumn = 15
$excel.ActiveWindow.ScrollRow = 1
$ws1.Range("W2").Select() | Out-Null

$ws2.Activate() | Out-Null
$excel.ActiveWindow.ScrollColumn = 5
$excel.ActiveWindow.ScrollRow = 1
$ws2.Range("E2").Select() | Out-Null
